# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Source data pulled from coinranking.com; row order for a few coins also
# changed between runs, so some rows get a full Coin/Link/Price/Volume swap
# rather than just updated numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ col = newValue }  (only cells that actually changed are listed)
$updates = @{
    2 = @{ "D" = '58.059.68' }
    3 = @{ "D" = '2.446.70'; "E" = '  -3.94%  ' }
    4 = @{ "E" = '  -0.04%  ' }
    5 = @{ "D" = '526.71'; "E" = '  -2.69%  ' }
    6 = @{ "D" = '133.07' }
    7 = @{ "D" = '0.999'; "E" = '  +0.37%  ' }
    8 = @{ "D" = '0.552'; "E" = '  -3.79%  ' }
    9 = @{ "D" = '2.449.77'; "E" = '  -4.69%  ' }
    10 = @{ "E" = '  -0.42%  ' }
    11 = @{ "D" = '0.0978'; "E" = '  -3.82%  ' }
    12 = @{ "D" = '5.30'; "E" = '  -3.30%  ' }
    13 = @{ "D" = '0.340'; "E" = '  -6.24%  ' }
    14 = @{ "D" = '2.881.33'; "E" = '  -3.81%  ' }
    15 = @{ "D" = '57.919.78'; "E" = '  -3.72%  ' }
    16 = @{ "D" = '22.43'; "E" = '  -7.75%  ' }
    17 = @{ "D" = '0.0000138'; "E" = '  -4.53%  ' }
    18 = @{ "D" = '2.451.74'; "E" = '  -5.13%  ' }
    19 = @{ "D" = '10.61'; "E" = '  -6.13%  ' }
    20 = @{ "D" = '318.63'; "E" = '  -2.89%  ' }
    21 = @{ "D" = '4.15'; "E" = '  -4.71%  ' }
    22 = @{ "E" = '  -0.22%  ' }
    23 = @{ "D" = '5.68'; "E" = '  -4.75%  ' }
    24 = @{ "D" = '62.10'; "E" = '  -1.61%  ' }
    25 = @{ "D" = '0.404'; "E" = '  -7.61%  ' }
    26 = @{ "D" = '0.164'; "E" = '  -2.24%  ' }
    27 = @{ "E" = '  -1.26%  ' }
    28 = @{ "D" = '7.41'; "E" = '  -7.70%  ' }
    29 = @{ "B" = 'Aptos'; "C" = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; "D" = '6.48'; "E" = '  -9.38%  ' }
    30 = @{ "B" = 'PEPE'; "C" = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; "D" = '0.0₃0744'; "E" = '  -7.30%  ' }
    31 = @{ "B" = 'PancakeSwap'; "C" = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; "D" = '1.74'; "E" = '  -4.44%  ' }
    32 = @{ "D" = '162.66'; "E" = '  -1.57%  ' }
    33 = @{ "E" = '  +0.10%  ' }
    34 = @{ "E" = '  -11.84%  ' }
    35 = @{ "B" = 'ImmutableX'; "C" = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; "D" = '1.35'; "E" = '  -9.31%  ' }
    36 = @{ "B" = 'EthereumClassic'; "C" = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; "D" = '18.09'; "E" = '  -3.80%  ' }
    37 = @{ "D" = '3.99'; "E" = '  -10.70%  ' }
    38 = @{ "E" = '  -7.49%  ' }
    39 = @{ "D" = '36.27'; "E" = '  -2.12%  ' }
    40 = @{ "D" = '3.50'; "E" = '  -6.53%  ' }
    41 = @{ "D" = '0.781'; "E" = '  -7.36%  ' }
    42 = @{ "E" = '  +0.36%  ' }
    43 = @{ "D" = '271.39'; "E" = '  -10.61%  ' }
    44 = @{ "D" = '5.00'; "E" = '  -11.05%  ' }
    45 = @{ "D" = '10.83'; "E" = '  -0.15%  ' }
    46 = @{ "D" = '0.584'; "E" = '  -4.41%  ' }
    47 = @{ "D" = '0.0916'; "E" = '  -2.65%  ' }
    48 = @{ "D" = '119.72'; "E" = '  -5.86%  ' }
    49 = @{ "D" = '0.0501'; "E" = '  -4.27%  ' }
    50 = @{ "D" = '0.0215'; "E" = '  -6.45%  ' }
    51 = @{ "D" = '16.80'; "E" = '  -8.47%  ' }
}

# Matches plain decimal numbers like "133.07" or "0.0978" (but not the
# thousand-dotted "58.059.68" style or the subscript-zero price strings,
# which already read back as text and need no special handling).
$numericLike = '^[+-]?[0-9]*\.?[0-9]+$'

foreach ($row in ($updates.Keys | Sort-Object)) {
    foreach ($col in $updates[$row].Keys) {
        $value = $updates[$row][$col]
        $cell = $ws.Range("$col$row")
        if ($col -eq "D" -and $value -match $numericLike) {
            # The Price column stores these as TEXT (e.g. "133.07"), but a bare
            # numeric-looking string assigned via .Value would be auto-coerced to
            # a real number by Excel. Use the quote-prefix convention (same as
            # typing '133.07 into the cell) to force text, then restore the Normal
            # style so no stray quote-prefix formatting is left behind.
            $cell.Value = "'" + $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
